# Updates cryptos list figures (price/volume) and reorders a few coin rows,
# matching the upstream GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D prices are stored as TEXT in the source file (e.g. "1.00", "65.960.32"),
# so every Column D write below is prefixed with a leading apostrophe to force
# Excel to keep it as text instead of normalising it into a Number.
$ws.Range("D2").Value = "'65.960.32"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "'3.301.14"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'573.78"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").Value = "'177.57"
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +4.41%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("D12").Value = "'3.876.33"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  -3.48%  "
$ws.Range("D14").Value = "'26.63"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").Value = "'66.070.52"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").Value = "'3.296.24"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "'437.10"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "'7.45"
$ws.Range("E21").Value = "  -4.31%  "
$ws.Range("D22").Value = "'72.62"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "'0.513"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'3.436.42"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  -4.08%  "
$ws.Range("E27").Value = "  +3.47%  "
$ws.Range("D28").Value = "'8.97"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("D31").Value = "'22.43"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'6.68"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'5.14"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("D35").Value = "'1.19"
$ws.Range("E35").Value = "  -3.56%  "
$ws.Range("E36").Value = "  -4.40%  "
$ws.Range("D37").Value = "'157.48"
$ws.Range("E37").Value = "  -3.22%  "
$ws.Range("D38").Value = "'27.01"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D40").Value = "'2.780.50"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").Value = "'0.785"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'6.14"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'40.37"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "'0.0660"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'321.86"
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.29"
$ws.Range("E47").Value = "  -4.80%  "
$ws.Range("D48").Value = "'23.54"
$ws.Range("E48").Value = "  -5.46%  "
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.13%  "

